# Apply rotation of observation records between rows 28-31 on the "Artfynd" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# --- Row 28 ---
$ws.Cells.Item(28, 1).Value = 131085484      # A28 Id
$ws.Cells.Item(28, 2).Value = 57884          # B28 Taxonsorteringsordning
$ws.Cells.Item(28, 5).Value = 100109         # E28 TaxonId
$ws.Cells.Item(28, 6).Value = "Tretåig hackspett"       # F28 Artnamn
$ws.Cells.Item(28, 7).Value = "Picoides tridactylus"    # G28 Vetenskapligt namn
$ws.Cells.Item(28, 8).Value = "(Linnaeus, 1758)"        # H28 Auktor
$ws.Cells.Item(28, 13).Value = "färska spår"            # M28 Aktivitet
$ws.Cells.Item(28, 17).Value = 585303        # Q28 Ost
$ws.Cells.Item(28, 18).Value = 7060488       # R28 Nord
$ws.Cells.Item(28, 19).Value = 15            # S28 Noggrannhet
$ws.Cells.Item(28, 26).Value = ""            # Z28 Starttid (removed)
$ws.Cells.Item(28, 28).Value = ""            # AB28 Sluttid (removed)
$ws.Cells.Item(28, 29).Value = "Färska ringhack gran"   # AC28 Publik kommentar
$ws.Cells.Item(28, 49).Value = "Daniel Rutschman"       # AW28 Rapportör
$ws.Cells.Item(28, 50).Value = "Daniel Rutschman"       # AX28 Observatörer

# --- Row 29 ---
$ws.Cells.Item(29, 1).Value = 131085240
$ws.Cells.Item(29, 2).Value = 57884
$ws.Cells.Item(29, 5).Value = 100109
$ws.Cells.Item(29, 6).Value = "Tretåig hackspett"
$ws.Cells.Item(29, 7).Value = "Picoides tridactylus"
$ws.Cells.Item(29, 8).Value = "(Linnaeus, 1758)"
$ws.Cells.Item(29, 13).Value = "färska spår"
$ws.Cells.Item(29, 17).Value = 585289
$ws.Cells.Item(29, 18).Value = 7060293
$ws.Cells.Item(29, 19).Value = 10
$ws.Cells.Item(29, 26).Value = "11:16"
$ws.Cells.Item(29, 28).Value = "11:16"
$ws.Cells.Item(29, 29).Value = "Ringhack på tall"
$ws.Cells.Item(29, 49).Value = "Kim Hultgren"
$ws.Cells.Item(29, 50).Value = "Kim Hultgren"

# --- Row 30 ---
$ws.Cells.Item(30, 1).Value = 131085178
$ws.Cells.Item(30, 2).Value = 91830
$ws.Cells.Item(30, 5).Value = 5432
$ws.Cells.Item(30, 6).Value = "Granticka"
$ws.Cells.Item(30, 7).Value = "Porodaedalea chrysoloma s.lat."
$ws.Cells.Item(30, 8).Value = ""
$ws.Cells.Item(30, 13).Value = ""            # M30 Aktivitet (removed)
$ws.Cells.Item(30, 17).Value = 585225
$ws.Cells.Item(30, 18).Value = 7060258
$ws.Cells.Item(30, 19).Value = 10
$ws.Cells.Item(30, 26).Value = "11:08"
$ws.Cells.Item(30, 28).Value = "11:08"
$ws.Cells.Item(30, 29).Value = ""            # AC30 Publik kommentar (removed)
$ws.Cells.Item(30, 49).Value = "Kim Hultgren"
$ws.Cells.Item(30, 50).Value = "Kim Hultgren"

# --- Row 31 ---
$ws.Cells.Item(31, 1).Value = 131085171
$ws.Cells.Item(31, 2).Value = 91806
$ws.Cells.Item(31, 5).Value = 1108
$ws.Cells.Item(31, 6).Value = "Harticka"
$ws.Cells.Item(31, 7).Value = "Pelloporus leporinus"
$ws.Cells.Item(31, 8).Value = "(Fr.) Krieglst."
$ws.Cells.Item(31, 13).Value = ""            # M31 Aktivitet (removed)
$ws.Cells.Item(31, 17).Value = 585222
$ws.Cells.Item(31, 18).Value = 7060254
$ws.Cells.Item(31, 19).Value = 15
$ws.Cells.Item(31, 26).Value = ""            # Z31 Starttid (removed)
$ws.Cells.Item(31, 28).Value = ""            # AB31 Sluttid (removed)
$ws.Cells.Item(31, 29).Value = ""            # AC31 Publik kommentar (removed)
$ws.Cells.Item(31, 49).Value = "Daniel Rutschman"
$ws.Cells.Item(31, 50).Value = "Daniel Rutschman"
